# The sheet holds weekly Haba (fava bean) price records for "Femacal de La
# Calera" in chronological order (by date, column D). A new weekly record was
# added ahead of the existing row that used to sit at row 66, pushing every
# subsequent record (old rows 66-103) down by one (to rows 67-104).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 66; this shifts the old rows 66-103 down to
# 67-104, growing the sheet's used range from A1:R103 to A1:R104.
$ws.Rows.Item(66).Insert()

# Fill in the newly inserted row 66 with the new weekly record.
$ws.Cells.Item(66, 1).Value = 3
$ws.Cells.Item(66, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(66, 3).Value = "Coquimbo"
$ws.Cells.Item(66, 4).Value = 44529
$ws.Cells.Item(66, 5).Value = 5
$ws.Cells.Item(66, 6).Value = 100112026
$ws.Cells.Item(66, 7).Value = "Haba"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 73
$ws.Cells.Item(66, 11).Value = 7000
$ws.Cells.Item(66, 12).Value = 7500
$ws.Cells.Item(66, 13).Value = 7240
$ws.Cells.Item(66, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(66, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(66, 16).Value = 290
$ws.Cells.Item(66, 17).Value = 25
$ws.Cells.Item(66, 18).Value = "Hortaliza"
